# Auto-generated edit script: updates market-board derived Leve profit data
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), as produced
# by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 975.36365
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 10500
$ws.Range("N38").Value = -11244
$ws.Range("H39").Value = 231.95454
$ws.Range("I39").Value = 93.8
$ws.Range("J39").Value = 528
$ws.Range("K39").Value = 281.4
$ws.Range("L39").Value = 1584
$ws.Range("M39").Value = 14.60000000000002
$ws.Range("N39").Value = -2176
$ws.Range("H62").Value = 7325.5
$ws.Range("I62").Value = 2651.25
$ws.Range("K62").Value = 2651.25
$ws.Range("M62").Value = -2027.25
$ws.Range("H65").Value = 7325.5
$ws.Range("I65").Value = 2651.25
$ws.Range("K65").Value = 13256.25
$ws.Range("M65").Value = -10136.25
$ws.Range("H80").Value = 407.13333
$ws.Range("I80").Value = 334.14285
$ws.Range("J80").Value = 471
$ws.Range("K80").Value = 1002.42855
$ws.Range("L80").Value = 1413
$ws.Range("M80").Value = -4.428550000000087
$ws.Range("N80").Value = -3409
$ws.Range("H83").Value = 407.13333
$ws.Range("I83").Value = 334.14285
$ws.Range("J83").Value = 471
$ws.Range("K83").Value = 3007.28565
$ws.Range("L83").Value = 4239
$ws.Range("M83").Value = 1984.71435
$ws.Range("N83").Value = -14223
$ws.Range("H137").Value = 2736.7307
$ws.Range("I137").Value = 1962.6666
$ws.Range("J137").Value = 3146.5293
$ws.Range("K137").Value = 5887.9998
$ws.Range("L137").Value = 9439.5879
$ws.Range("M137").Value = -3337.9998
$ws.Range("N137").Value = -14539.5879

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1183.875
$ws.Range("I2").Value = 1183.875
$ws.Range("K2").Value = 1183.875
$ws.Range("M2").Value = -1070.875
$ws.Range("H32").Value = 4169022.5
$ws.Range("I32").Value = 577.15
$ws.Range("K32").Value = 577.15
$ws.Range("M32").Value = -290.15
$ws.Range("H45").Value = 2477.5454
$ws.Range("J45").Value = 2843.5
$ws.Range("L45").Value = 2843.5
$ws.Range("N45").Value = -3597.5
$ws.Range("H61").Value = 1994.0526
$ws.Range("I61").Value = 1994.0526
$ws.Range("K61").Value = 1994.0526
$ws.Range("M61").Value = -1782.0526
$ws.Range("H116").Value = 1183.875
$ws.Range("I116").Value = 1183.875
$ws.Range("K116").Value = 1183.875
$ws.Range("M116").Value = 1110.125
$ws.Range("H122").Value = 1466.409
$ws.Range("I122").Value = 1273.4117
$ws.Range("K122").Value = 3820.2351
$ws.Range("M122").Value = -1370.2351
$ws.Range("H136").Value = 1994.0526
$ws.Range("I136").Value = 1994.0526
$ws.Range("K136").Value = 5982.1578
$ws.Range("M136").Value = -3432.1578

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1183.875
$ws.Range("I3").Value = 1183.875
$ws.Range("K3").Value = 1183.875
$ws.Range("M3").Value = -1069.875
$ws.Range("H105").Value = 1170
$ws.Range("I105").Value = 1170
$ws.Range("K105").Value = 1170
$ws.Range("M105").Value = 577
$ws.Range("H134").Value = 4027.45
$ws.Range("I134").Value = 997.05884
$ws.Range("J134").Value = 21199.666
$ws.Range("K134").Value = 2991.17652
$ws.Range("L134").Value = 63598.99800000001
$ws.Range("M134").Value = -456.17652
$ws.Range("N134").Value = -68668.99800000001

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2991.5
$ws.Range("I58").Value = 1790.8
$ws.Range("K58").Value = 1790.8
$ws.Range("M58").Value = -1587.8
$ws.Range("H99").Value = 2363.4443
$ws.Range("I99").Value = 1536.7142
$ws.Range("K99").Value = 1536.7142
$ws.Range("M99").Value = -38.71419999999989
$ws.Range("H102").Value = 28160.666
$ws.Range("J102").Value = 28160.666
$ws.Range("L102").Value = 28160.666
$ws.Range("N102").Value = -33028.666
$ws.Range("H122").Value = 592.2727
$ws.Range("I122").Value = 758.7143
$ws.Range("J122").Value = 301
$ws.Range("K122").Value = 2276.1429
$ws.Range("L122").Value = 903
$ws.Range("M122").Value = 173.8571000000002
$ws.Range("N122").Value = -5803
$ws.Range("H126").Value = 2363.4443
$ws.Range("I126").Value = 1536.7142
$ws.Range("K126").Value = 4610.142599999999
$ws.Range("M126").Value = -2140.142599999999
$ws.Range("H132").Value = 1828.2759
$ws.Range("I132").Value = 1462.7307
$ws.Range("J132").Value = 4996.3335
$ws.Range("K132").Value = 4388.1921
$ws.Range("L132").Value = 14989.0005
$ws.Range("M132").Value = -1858.1921
$ws.Range("N132").Value = -20049.0005
$ws.Range("H136").Value = 2991.5
$ws.Range("I136").Value = 1790.8
$ws.Range("K136").Value = 5372.4
$ws.Range("M136").Value = -2822.4

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2395.3845
$ws.Range("I5").Value = 1822.2222
$ws.Range("K5").Value = 5466.6666
$ws.Range("M5").Value = -5354.6666
$ws.Range("H128").Value = 428567
$ws.Range("I128").Value = 428567
$ws.Range("K128").Value = 1285701
$ws.Range("M128").Value = -1280721
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -50060
$ws.Range("H135").Value = 2395.3845
$ws.Range("I135").Value = 1822.2222
$ws.Range("K135").Value = 16399.9998
$ws.Range("M135").Value = -13864.9998

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20009
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 20009
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 4874.75
$ws.Range("I80").Value = 4833
$ws.Range("K80").Value = 4833
$ws.Range("M80").Value = -3835
$ws.Range("H83").Value = 4874.75
$ws.Range("I83").Value = 4833
$ws.Range("K83").Value = 24165
$ws.Range("M83").Value = -19173
$ws.Range("H122").Value = 1154
$ws.Range("I122").Value = 604
$ws.Range("J122").Value = 2254
$ws.Range("K122").Value = 1812
$ws.Range("L122").Value = 6762
$ws.Range("M122").Value = 638
$ws.Range("N122").Value = -11662
$ws.Range("H132").Value = 66352.44
$ws.Range("I132").Value = 86571.164
$ws.Range("J132").Value = 5696.25
$ws.Range("K132").Value = 259713.492
$ws.Range("L132").Value = 17088.75
$ws.Range("M132").Value = -257183.492
$ws.Range("N132").Value = -22148.75

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3319
$ws.Range("I61").Value = 1891
$ws.Range("J61").Value = 6175
$ws.Range("K61").Value = 1891
$ws.Range("L61").Value = 6175
$ws.Range("M61").Value = -1689
$ws.Range("N61").Value = -6579
$ws.Range("H102").Value = 70561
$ws.Range("J102").Value = 70561
$ws.Range("L102").Value = 70561
$ws.Range("N102").Value = -77051
$ws.Range("H113").Value = 3319
$ws.Range("I113").Value = 1891
$ws.Range("J113").Value = 6175
$ws.Range("K113").Value = 1891
$ws.Range("L113").Value = 6175
$ws.Range("M113").Value = 279
$ws.Range("N113").Value = -10515
$ws.Range("H132").Value = 5150.6665
$ws.Range("I132").Value = 5150.6665
$ws.Range("K132").Value = 15451.9995
$ws.Range("M132").Value = -12921.9995
$ws.Range("H136").Value = 2859.2
$ws.Range("I136").Value = 2823.375
$ws.Range("K136").Value = 8470.125
$ws.Range("M136").Value = -5920.125
$ws.Range("J140").Value = 200000
$ws.Range("L140").Value = 200000
$ws.Range("N140").Value = -210360

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 70000
$ws.Range("J86").Value = 70000
$ws.Range("L86").Value = 70000
$ws.Range("N86").Value = -72246
$ws.Range("H89").Value = 70000
$ws.Range("J89").Value = 70000
$ws.Range("L89").Value = 350000
$ws.Range("N89").Value = -361232
$ws.Range("H99").Value = 60000
$ws.Range("I99").Value = 60000
$ws.Range("K99").Value = 60000
$ws.Range("M99").Value = -57005
$ws.Range("H102").Value = 100336.5
$ws.Range("J102").Value = 100336.5
$ws.Range("L102").Value = 100336.5
$ws.Range("N102").Value = -106826.5
$ws.Range("H106").Value = 58888
$ws.Range("J106").Value = 58888
$ws.Range("L106").Value = 58888
$ws.Range("N106").Value = -61412
$ws.Range("H113").Value = 1087.2142
$ws.Range("I113").Value = 953.3
$ws.Range("J113").Value = 1422
$ws.Range("K113").Value = 2859.9
$ws.Range("L113").Value = 4266
$ws.Range("M113").Value = -689.8999999999996
$ws.Range("N113").Value = -8606
